# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest count) figures and flip two events
# that have sold out (最低票价 -> "不可售") across all four sheets.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 26910
$ws.Range("F4").Value  = 600
$ws.Range("G4").Value  = "不可售"
$ws.Range("F6").Value  = 623
$ws.Range("F7").Value  = 180
$ws.Range("F10").Value = 367
$ws.Range("F13").Value = 51
$ws.Range("F15").Value = 85
$ws.Range("F16").Value = 451
$ws.Range("F18").Value = 1579
$ws.Range("F19").Value = 225
$ws.Range("F20").Value = 63

# --- 演出 (Performances) sheet ----------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 38
$ws.Range("F10").Value = 443
$ws.Range("F12").Value = 8
$ws.Range("F14").Value = 18

# --- 本地生活 (Local life) sheet --------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5125
$ws.Range("F3").Value = 251

# --- 全部类型 (All types) sheet, mirrors the rows above ---------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 5125
$ws.Range("F4").Value  = 251
$ws.Range("F5").Value  = 26910
$ws.Range("F6").Value  = 600
$ws.Range("G6").Value  = "不可售"
$ws.Range("F10").Value = 623
$ws.Range("F13").Value = 180
$ws.Range("F16").Value = 38
$ws.Range("F18").Value = 443
$ws.Range("F22").Value = 367
$ws.Range("F25").Value = 51
$ws.Range("F26").Value = 8
$ws.Range("F28").Value = 85
$ws.Range("F30").Value = 18
$ws.Range("F31").Value = 451
$ws.Range("F34").Value = 1579
$ws.Range("F35").Value = 225
$ws.Range("F37").Value = 63
